{"js": "// The second paragraph of the document contains a single Word field whose\n// instruction text is:  m:'doc.html'.fromHTMLURI()\n// The edit turns that field into plain literal text reading\n// \"{m:'doc.html'.fromHTMLURI()}\" (braces added, field machinery removed),\n// while keeping the \"_GoBack\" bookmark that sits in the middle of the text\n// (right after \"doc.html\").\n\nconst body = context.document.body;\n\n// Locate the field (there is exactly one in this document).\nconst fields = body.fields;\nfields.load(\"items\");\nawait context.sync();\n\nconst field = fields.items[0];\n\n// Turn the field into its (empty) result - this removes the fldChar/\n// instrText plumbing but keeps the owning paragraph's identity/properties\n// intact (it becomes an empty paragraph).\nfield.unlink();\nawait context.sync();\n\n// Re-fetch paragraphs and grab the (now empty) paragraph that used to hold\n// the field.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst fieldParagraph = paragraphs.items[1];\n\n// Rebuild the literal text, preserving the bookmark in the same spot it\n// occupied inside the original field instruction (right after \"doc.html\").\nconst ooxml = `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">` +\n  `<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">` +\n  `<pkg:xmlData>` +\n  `<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">` +\n  `<w:body>` +\n  `<w:p>` +\n  `<w:r><w:t>{</w:t></w:r>` +\n  `<w:r><w:t>m</w:t></w:r>` +\n  `<w:r><w:t>:</w:t></w:r>` +\n  `<w:r><w:t>'</w:t></w:r>` +\n  `<w:r><w:t>doc.html</w:t></w:r>` +\n  `<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>` +\n  `<w:bookmarkEnd w:id=\"0\"/>` +\n  `<w:r><w:t>'.fromHTMLURI()</w:t></w:r>` +\n  `<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>` +\n  `</w:p>` +\n  `</w:body>` +\n  `</w:document>` +\n  `</pkg:xmlData>` +\n  `</pkg:part>` +\n  `</pkg:package>`;\n\nfieldParagraph.insertOoxml(ooxml, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# The second paragraph of the document holds a single Word field whose\n# instruction text is:  m:'doc.html'.fromHTMLURI()\n# The edit turns that field into plain literal text reading\n# \"{m:'doc.html'.fromHTMLURI()}\" (braces added, field machinery removed),\n# while keeping the \"_GoBack\" bookmark exactly where it sat inside the\n# original field instruction text (right after \"doc.html\").\n\n$d = $word.ActiveDocument\n\n# There is exactly one field in the document - the one holding the query\n# instruction text - and it lives alone in the 2nd paragraph.\n$field = $d.Fields.Item(1)\n\n# Remember where that paragraph starts before we touch anything (Unlink\n# collapses the field down to its - empty - result; the paragraph itself,\n# along with its rsid/formatting attributes, is left in place and keeps\n# the same start offset).\n$paragraphStart = $d.Paragraphs(2).Range.Start\n\n$field.Unlink()\n\n# Re-resolve the (now empty) paragraph and type the literal text that used\n# to be the field code, wrapped in braces.\n$p = $d.Paragraphs(2)\n$literalText = \"{m:'doc.html'.fromHTMLURI()}\"\n$p.Range.InsertAfter($literalText)\n\n# Figure out the character offsets of each token boundary so we can force\n# Word to split the text into one run per token, the same way the field's\n# instruction text used to be split into one run per instrText token.\n$tokens = @(\"{\", \"m\", \":\", \"'\", \"doc.html\", \"'.fromHTMLURI()\", \"}\")\n$boundaries = @()\n$offset = $paragraphStart\nfor ($i = 0; $i -lt $tokens.Length - 1; $i++) {\n    $offset += $tokens[$i].Length\n    $boundaries += $offset\n}\n\n# The bookmark used to sit right after \"doc.html\" (4th token, 0-based index\n# 4 once incremented) inside the field instruction text - keep it there.\n$goBackBoundary = $boundaries[4]\n\n# Temporary bookmarks force a run split at each token boundary (Word never\n# merges runs across a bookmark). The permanent \"_GoBack\" bookmark is added\n# directly at its final position; the rest are removed once they have done\n# their job, leaving plain run boundaries behind.\nfor ($i = 0; $i -lt $boundaries.Length; $i++) {\n    $pos = $boundaries[$i]\n    $splitRange = $d.Range($pos, $pos)\n    if ($pos -eq $goBackBoundary) {\n        $d.Bookmarks.Add(\"_GoBack\", $splitRange)\n    } else {\n        $d.Bookmarks.Add(\"TempSplit$i\", $splitRange)\n    }\n}\nfor ($i = 0; $i -lt $boundaries.Length; $i++) {\n    if ($boundaries[$i] -ne $goBackBoundary) {\n        $d.Bookmarks(\"TempSplit$i\").Delete()\n    }\n}\n"}
